$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# =========================================================
# Step 1: "Parameters" sheet -- add new SAVINGS_RATE row (37)
# =========================================================
$ws1.Range("A37").Value = "SAVINGS_RATE"
$ws1.Range("B37").Value = 0.056
$ws1.Range("A37:B37").Font.Name = "Aptos Narrow"
$ws1.Range("A37:B37").Font.Size = 12
$ws1.Range("A37:B37").Font.Bold = $false
$ws1.Rows.Item(37).RowHeight = 16

# =========================================================
# Step 2: "Parameters" sheet -- un-bold column A, rows 23-36
#   (previously bold Helvetica Neue, now matches column B)
# =========================================================
$ws1.Range("A23:B36").Font.Name = "Helvetica Neue"
$ws1.Range("A23:B36").Font.Size = 10
$ws1.Range("A23:B36").Font.Bold = $false

# =========================================================
# Step 3: add the new "Info" sheet right after "Parameters"
# =========================================================
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Info"

$newSheet.Columns.Item(1).ColumnWidth = 58.33203125
$newSheet.Columns.Item(2).ColumnWidth = 58.83203125

# --- header row 3: "Parameter Name" / "Description" (bold Aptos Narrow 12) ---
$newSheet.Range("A3").Value = "Parameter Name"
$newSheet.Range("B3").Value = "Description"
$newSheet.Range("A3:B3").Font.Name = "Aptos Narrow"
$newSheet.Range("A3:B3").Font.Size = 12
$newSheet.Range("A3:B3").Font.Bold = $true
$newSheet.Rows.Item(3).RowHeight = 16

# --- rows 4-38: KEY name (col A, Arial Unicode MS 10) + description (col B, Aptos Narrow 12) ---
$newSheet.Range("A4").Value = "MIN_AGE_TO_HAVE_INCOME"
$newSheet.Range("B4").Value = "Minimum age to have non-employment, non-benefit income"
$newSheet.Range("A5").Value = "MAX_LABOUR_HOURS_IN_WEEK"
$newSheet.Range("B5").Value = "Maximum number of labour hours allowed in a week"
$newSheet.Range("A6").Value = "HOURS_IN_WEEK"
$newSheet.Range("B6").Value = "Total hours in a week (used to calculate leisure in labour supply)"
$newSheet.Range("A7").Value = "USE_CONTINUOUS_LABOUR_SUPPLY_HOURS"
$newSheet.Range("B7").Value = "If true, generates random weekly labour supply hours within each bracket; if false, uses fixed hours for all persons"
$newSheet.Range("A8").Value = "AGE_TO_BECOME_RESPONSIBLE"
$newSheet.Range("B8").Value = "Age at which a person becomes reference person of their own benefit unit"
$newSheet.Range("A9").Value = "MIN_AGE_TO_LEAVE_EDUCATION"
$newSheet.Range("B9").Value = "Minimum age to leave full-time education"
$newSheet.Range("A10").Value = "MAX_AGE_TO_LEAVE_CONTINUOUS_EDUCATION"
$newSheet.Range("B10").Value = "Maximum age to remain in continuous education"
$newSheet.Range("A11").Value = "MAX_AGE_TO_ENTER_EDUCATION"
$newSheet.Range("B11").Value = "Maximum age to enter education"
$newSheet.Range("A12").Value = "MIN_AGE_TO_RETIRE"
$newSheet.Range("B12").Value = "Minimum age to consider retirement"
$newSheet.Range("A13").Value = "DEFAULT_AGE_TO_RETIRE"
$newSheet.Range("B13").Value = "Default retirement age (if pension included but retirement decision not modeled)"
$newSheet.Range("A14").Value = "MIN_AGE_FORMAL_SOCARE"
$newSheet.Range("B14").Value = "Minimum age to receive formal social care"
$newSheet.Range("A15").Value = "MIN_AGE_FLEXIBLE_LABOUR_SUPPLY"
$newSheet.Range("B15").Value = "Minimum age for flexible labour supply eligibility"
$newSheet.Range("A16").Value = "MAX_AGE_FLEXIBLE_LABOUR_SUPPLY"
$newSheet.Range("B16").Value = "Maximum age for flexible labour supply eligibility"
$newSheet.Range("A17").Value = "SHARE_OF_WEALTH_TO_ANNUITISE_AT_RETIREMENT"
$newSheet.Range("B17").Value = "Proportion of wealth to annuitise at retirement"
$newSheet.Range("A18").Value = "ANNUITY_RATE_OF_RETURN"
$newSheet.Range("B18").Value = "Assumed annuity rate of return"
$newSheet.Range("A19").Value = "MIN_HOURS_FULL_TIME_EMPLOYED"
$newSheet.Range("B19").Value = "Minimum weekly hours defining full-time employment"
$newSheet.Range("A20").Value = "MIN_HOURLY_WAGE_RATE"
$newSheet.Range("B20").Value = "Minimum possible hourly wage"
$newSheet.Range("A21").Value = "MAX_HOURLY_WAGE_RATE"
$newSheet.Range("B21").Value = "Maximum possible hourly wage"
$newSheet.Range("A22").Value = "MAX_HOURS_WEEKLY_FORMAL_CARE"
$newSheet.Range("B22").Value = "Maximum number of hours of formal care per week"
$newSheet.Range("A23").Value = "MAX_HOURS_WEEKLY_INFORMAL_CARE"
$newSheet.Range("B23").Value = "Maximum number of hours of informal care per week"
$newSheet.Range("A24").Value = "CHILDCARE_COST_EARNINGS_CAP"
$newSheet.Range("B24").Value = "Maximum share of earnings payable as childcare"
$newSheet.Range("A25").Value = "TAXDB_REGIMES"
$newSheet.Range("B25").Value = "Number of tax/benefit regimes supported"
$newSheet.Range("A26").Value = "MIN_START_YEAR"
$newSheet.Range("B26").Value = "Minimum allowed simulation start year (oldest initial population)"
$newSheet.Range("A27").Value = "MAX_START_YEAR"
$newSheet.Range("B27").Value = "Maximum allowed simulation start year (most recent initial population)"
$newSheet.Range("A28").Value = "MIN_START_YEAR_TRAINING"
$newSheet.Range("B28").Value = "Minimum allowed training start year"
$newSheet.Range("A29").Value = "MAX_START_YEAR_TRAINING"
$newSheet.Range("B29").Value = "Maximum allowed training start year"
$newSheet.Range("A30").Value = "MIN_CAPITAL_INCOME_PER_MONTH"
$newSheet.Range("B30").Value = "Minimum capital income per month"
$newSheet.Range("A31").Value = "MAX_CAPITAL_INCOME_PER_MONTH"
$newSheet.Range("B31").Value = "Maximum capital income per month"
$newSheet.Range("A32").Value = "MIN_PERSONAL_PENSION_PER_MONTH"
$newSheet.Range("B32").Value = "Minimum pension income per month"
$newSheet.Range("A33").Value = "MAX_PERSONAL_PENSION_PER_MONTH"
$newSheet.Range("B33").Value = "Maximum pension income per month"
$newSheet.Range("A34").Value = "MAX_CHILD_AGE_FOR_FORMAL_CARE"
$newSheet.Range("B34").Value = "Maximum age of child eligible for formal care"
$newSheet.Range("A35").Value = "MIN_AGE_MATERNITY"
$newSheet.Range("B35").Value = "Minimum age a person can give birth"
$newSheet.Range("A36").Value = "MAX_AGE_MATERNITY"
$newSheet.Range("B36").Value = "Maximum age a person can give birth"
$newSheet.Range("A37").Value = "BASE_PRICE_YEAR"
$newSheet.Range("B37").Value = "Base year for model parameters (prices)"
$newSheet.Range("A38").Value = "PROB_NEWBORN_IS_MALE"
$newSheet.Range("B38").Value = "Probability a newborn is male"

$dataRange = $newSheet.Range("A4:B38")
$dataRange.Rows.RowHeight = 17

$newSheet.Range("A4:A38").Font.Name = "Arial Unicode MS"
$newSheet.Range("A4:A38").Font.Size = 10
$newSheet.Range("A4:A38").Font.Bold = $false

$newSheet.Range("B4:B38").Font.Name = "Aptos Narrow"
$newSheet.Range("B4:B38").Font.Size = 12
$newSheet.Range("B4:B38").Font.Bold = $false

# --- final row: SAVINGS_RATE entry (same style as rows 4-38 but ht=16, like header/footer rows) ---
$newSheet.Range("A39").Value = "SAVINGS_RATE"
$newSheet.Range("B39").Value = "Country-specific savings rate"
$newSheet.Range("A39:B39").Font.Name = "Aptos Narrow"
$newSheet.Range("A39:B39").Font.Size = 12
$newSheet.Range("A39:B39").Font.Bold = $false
$newSheet.Rows.Item(39).RowHeight = 16

# --- row 2: blank spacer row (same font as surrounding) ---
$newSheet.Range("A2:B2").Font.Name = "Aptos Narrow"
$newSheet.Range("A2:B2").Font.Size = 12
$newSheet.Range("A2:B2").Font.Bold = $false
$newSheet.Rows.Item(2).RowHeight = 16

# --- row 1: intro sentence (set LAST so its string lands at the end of sharedStrings) ---
$newSheet.Range("A1").Value = "This file is used to rewrite the following country-specific parameters"
$newSheet.Range("A1:B1").Font.Name = "Aptos Narrow"
$newSheet.Range("A1:B1").Font.Size = 12
$newSheet.Range("A1:B1").Font.Bold = $false
$newSheet.Rows.Item(1).RowHeight = 16

# =========================================================
# Step 4: selections -- Parameters -> A12 ; Info -> A1:B39 (active tab)
# =========================================================
$ws1.Range("A12").Select()
$newSheet.Range("A1:B39").Select()

